$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H6: ratio of cost to s_nom for the 380kV line
$ws.Range("H6").Formula = "=G6/C6"
$ws.Range("H6").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# New headers for parallel cables / cost computation in row 8
$ws.Range("L8").Value = "parallel cables"
$ws.Range("M8").Value = "cost"

# Number of parallel cables (18 circuits / 3 per something = 6)
$ws.Range("L9").Formula = "=18/3"

# Cost per km per MW computed from distance, resistance ratio and parallel cables
$ws.Range("M9").Formula = "=G6/J9*K9/L9"
$ws.Range("M9").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# Fit the newly used columns to their content (matches Excel's bestFit behaviour
# as closely as this engine's column-width model allows)
$ws.Columns.Item(8).ColumnWidth = 8.333333333333334
$ws.Columns.Item(12).ColumnWidth = 11.833333333333334
$ws.Columns.Item(13).ColumnWidth = 10.333333333333334

# Update active selection to reflect where the user ended up after editing
$ws.Range("O11").Select()
